$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 397; $r++) {
    $ws.Cells.Item($r, 3).Value = 45177
}
